# Update crypto price/volume figures per the latest GitHub Actions refresh.
# Values are written with a leading apostrophe so Excel stores them as literal
# text (preserving exact digits, trailing zeros, and the padded percent
# strings) instead of silently re-parsing them as floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'28.552.24"
$ws.Range("E2").Value = "'  +1.28%  "

# Row 3
$ws.Range("D3").Value = "'1.912.01"
$ws.Range("E3").Value = "'  +4.67%  "

# Row 4
$ws.Range("E4").Value = "'  +0.18%  "

# Row 5
$ws.Range("D5").Value = "'315.16"
$ws.Range("E5").Value = "'  +1.35%  "

# Row 6
$ws.Range("E6").Value = "'  +0.07%  "

# Row 7
$ws.Range("D7").Value = "'0.5158"
$ws.Range("E7").Value = "'  +3.88%  "

# Row 8
$ws.Range("D8").Value = "'0.3954"
$ws.Range("E8").Value = "'  +0.65%  "

# Row 9
$ws.Range("D9").Value = "'0.09660"
$ws.Range("E9").Value = "'  -2.45%  "

# Row 10
$ws.Range("E10").Value = "'  +3.70%  "

# Row 11
$ws.Range("D11").Value = "'42.03"
$ws.Range("E11").Value = "'  +1.80%  "

# Row 12
$ws.Range("D12").Value = "'6.529"
$ws.Range("E12").Value = "'  +1.19%  "

# Row 13
$ws.Range("D13").Value = "'21.24"
$ws.Range("E13").Value = "'  +2.83%  "

# Row 14
$ws.Range("D14").Value = "'1.917.21"
$ws.Range("E14").Value = "'  +5.25%  "

# Row 15
$ws.Range("D15").Value = "'7.502"
$ws.Range("E15").Value = "'  +2.62%  "

# Row 16
$ws.Range("E16").Value = "'  +0.15%  "

# Row 17
$ws.Range("D17").Value = "'94.63"
$ws.Range("E17").Value = "'  +1.99%  "

# Row 18
$ws.Range("D18").Value = "'0.00001133"
$ws.Range("E18").Value = "'  -0.96%  "

# Row 19
$ws.Range("D19").Value = "'0.06657"
$ws.Range("E19").Value = "'  -0.07%  "

# Row 20
$ws.Range("D20").Value = "'18.21"
$ws.Range("E20").Value = "'  +5.55%  "

# Row 21
$ws.Range("E21").Value = "'  +0.00%  "

# Row 22
$ws.Range("D22").Value = "'6.313"
$ws.Range("E22").Value = "'  +5.26%  "

# Row 23
$ws.Range("D23").Value = "'28.613.80"
$ws.Range("E23").Value = "'  +1.36%  "

# Row 24
$ws.Range("D24").Value = "'11.51"
$ws.Range("E24").Value = "'  +1.29%  "

# Row 25
$ws.Range("D25").Value = "'2.311"
$ws.Range("E25").Value = "'  +2.99%  "

# Row 26
$ws.Range("D26").Value = "'2.677"
$ws.Range("E26").Value = "'  +10.46%  "

# Row 27
$ws.Range("D27").Value = "'2.134.15"
$ws.Range("E27").Value = "'  +5.01%  "

# Row 29
$ws.Range("D29").Value = "'158.33"
$ws.Range("E29").Value = "'  -0.30%  "

# Row 30
$ws.Range("D30").Value = "'128.73"
$ws.Range("E30").Value = "'  +1.20%  "

# Row 31
$ws.Range("D31").Value = "'1.110"
$ws.Range("E31").Value = "'  +6.51%  "

# Row 32
$ws.Range("D32").Value = "'0.1078"
$ws.Range("E32").Value = "'  +2.08%  "

# Row 33
$ws.Range("D33").Value = "'5.760"
$ws.Range("E33").Value = "'  +2.62%  "

# Row 34
$ws.Range("D34").Value = "'3.635"
$ws.Range("E34").Value = "'  +0.65%  "

# Row 35
$ws.Range("D35").Value = "'10.11"
$ws.Range("E35").Value = "'  +11.81%  "

# Row 36
$ws.Range("D36").Value = "'0.06785"
$ws.Range("E36").Value = "'  +0.44%  "

# Row 37
$ws.Range("D37").Value = "'1.278"
$ws.Range("E37").Value = "'  +8.18%  "

# Row 38
$ws.Range("D38").Value = "'0.02432"
$ws.Range("E38").Value = "'  +3.72%  "

# Row 39
$ws.Range("D39").Value = "'0.2216"
$ws.Range("E39").Value = "'  +2.69%  "

# Row 40
$ws.Range("D40").Value = "'11.81"
$ws.Range("E40").Value = "'  +3.47%  "

# Row 41
$ws.Range("D41").Value = "'5.085"
$ws.Range("E41").Value = "'  +1.95%  "

# Row 42
$ws.Range("D42").Value = "'0.6461"
$ws.Range("E42").Value = "'  +3.65%  "

# Row 43
$ws.Range("E43").Value = "'  +0.66%  "

# Row 44
$ws.Range("E44").Value = "'  -0.01%  "

# Row 45
$ws.Range("D45").Value = "'13.52"
$ws.Range("E45").Value = "'  +1.96%  "

# Row 46
$ws.Range("D46").Value = "'0.6103"
$ws.Range("E46").Value = "'  +2.56%  "

# Row 47
$ws.Range("D47").Value = "'3.779"
$ws.Range("E47").Value = "'  +1.84%  "

# Row 48
$ws.Range("D48").Value = "'1.283"
$ws.Range("E48").Value = "'  +0.86%  "

# Row 49
$ws.Range("D49").Value = "'2.034"
$ws.Range("E49").Value = "'  +4.24%  "

# Row 50
$ws.Range("D50").Value = "'125.10"
$ws.Range("E50").Value = "'  +0.68%  "

# Row 51
$ws.Range("E51").Value = "'  +1.55%  "
